# Auto-generated edit script: numeric corrections to Leve profit calcs
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 7936.143  # H62: 7901.3335 -> 7936.143
$ws.Cells.Item(62, 10).Value = 7961.1665  # J62: 7915.75 -> 7961.1665
$ws.Cells.Item(62, 12).Value = 7961.1665  # L62: 7915.75 -> 7961.1665
$ws.Cells.Item(62, 14).Value = -9209.166499999999  # N62: -9163.75 -> -9209.166499999999
$ws.Cells.Item(65, 8).Value = 7936.143  # H65: 7901.3335 -> 7936.143
$ws.Cells.Item(65, 10).Value = 7961.1665  # J65: 7915.75 -> 7961.1665
$ws.Cells.Item(65, 12).Value = 39805.8325  # L65: 39578.75 -> 39805.8325
$ws.Cells.Item(65, 14).Value = -46045.8325  # N65: -45818.75 -> -46045.8325
$ws.Cells.Item(86, 8).Value = 3500.353  # H86: 3663.5625 -> 3500.353
$ws.Cells.Item(86, 9).Value = 1793.7273  # I86: 1926.8889 -> 1793.7273
$ws.Cells.Item(86, 10).Value = 6629.1665  # J86: 5896.4287 -> 6629.1665
$ws.Cells.Item(86, 11).Value = 1793.7273  # K86: 1926.8889 -> 1793.7273
$ws.Cells.Item(86, 12).Value = 6629.1665  # L86: 5896.4287 -> 6629.1665
$ws.Cells.Item(86, 13).Value = -670.7273  # M86: -803.8888999999999 -> -670.7273
$ws.Cells.Item(86, 14).Value = -8875.166499999999  # N86: -8142.4287 -> -8875.166499999999
$ws.Cells.Item(89, 8).Value = 3500.353  # H89: 3663.5625 -> 3500.353
$ws.Cells.Item(89, 9).Value = 1793.7273  # I89: 1926.8889 -> 1793.7273
$ws.Cells.Item(89, 10).Value = 6629.1665  # J89: 5896.4287 -> 6629.1665
$ws.Cells.Item(89, 11).Value = 8968.636500000001  # K89: 9634.4445 -> 8968.636500000001
$ws.Cells.Item(89, 12).Value = 33145.8325  # L89: 29482.1435 -> 33145.8325
$ws.Cells.Item(89, 13).Value = -3352.636500000001  # M89: -4018.4445 -> -3352.636500000001
$ws.Cells.Item(89, 14).Value = -44377.8325  # N89: -40714.14350000001 -> -44377.8325
$ws.Cells.Item(106, 8).Value = 926.6667  # H106: 965 -> 926.6667
$ws.Cells.Item(106, 9).Value = 926.6667  # I106: 965 -> 926.6667
$ws.Cells.Item(106, 11).Value = 926.6667  # K106: 965 -> 926.6667
$ws.Cells.Item(106, 13).Value = -295.6667  # M106: -334 -> -295.6667
$ws.Cells.Item(132, 8).Value = 1069.5555  # H132: 1049.579 -> 1069.5555
$ws.Cells.Item(132, 9).Value = 1015.9375  # I132: 996.7646999999999 -> 1015.9375
$ws.Cells.Item(132, 11).Value = 3047.8125  # K132: 2990.2941 -> 3047.8125
$ws.Cells.Item(132, 13).Value = -517.8125  # M132: -460.2941000000001 -> -517.8125
$ws.Cells.Item(137, 8).Value = 1706.1818  # H137: 1735.6666 -> 1706.1818
$ws.Cells.Item(137, 10).Value = 2070  # J137: 2065 -> 2070
$ws.Cells.Item(137, 12).Value = 6210  # L137: 6195 -> 6210
$ws.Cells.Item(137, 14).Value = -11310  # N137: -11295 -> -11310

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1150.2727  # H2: 1245.3 -> 1150.2727
$ws.Cells.Item(2, 9).Value = 572.55554  # I2: 619.125 -> 572.55554
$ws.Cells.Item(2, 11).Value = 572.55554  # K2: 619.125 -> 572.55554
$ws.Cells.Item(2, 13).Value = -459.55554  # M2: -506.125 -> -459.55554
$ws.Cells.Item(45, 8).Value = 802.5  # H45: 1036.6666 -> 802.5
$ws.Cells.Item(45, 9).Value = 802.5  # I45: 1036.6666 -> 802.5
$ws.Cells.Item(45, 11).Value = 802.5  # K45: 1036.6666 -> 802.5
$ws.Cells.Item(45, 13).Value = -425.5  # M45: -659.6666 -> -425.5
$ws.Cells.Item(97, 8).Value = 796.3333  # H97: 560.6667 -> 796.3333
$ws.Cells.Item(97, 9).Value = 755.6  # I97: 505.75 -> 755.6
$ws.Cells.Item(97, 11).Value = 755.6  # K97: 505.75 -> 755.6
$ws.Cells.Item(97, 13).Value = -259.6  # M97: -9.75 -> -259.6
$ws.Cells.Item(116, 8).Value = 1150.2727  # H116: 1245.3 -> 1150.2727
$ws.Cells.Item(116, 9).Value = 572.55554  # I116: 619.125 -> 572.55554
$ws.Cells.Item(116, 11).Value = 572.55554  # K116: 619.125 -> 572.55554
$ws.Cells.Item(116, 13).Value = 1721.44446  # M116: 1674.875 -> 1721.44446

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1150.2727  # H3: 1245.3 -> 1150.2727
$ws.Cells.Item(3, 9).Value = 572.55554  # I3: 619.125 -> 572.55554
$ws.Cells.Item(3, 11).Value = 572.55554  # K3: 619.125 -> 572.55554
$ws.Cells.Item(3, 13).Value = -458.55554  # M3: -505.125 -> -458.55554
$ws.Cells.Item(6, 8).Value = 30499.5  # H6: 30994 -> 30499.5
$ws.Cells.Item(6, 10).Value = 30499.5  # J6: 30994 -> 30499.5
$ws.Cells.Item(6, 12).Value = 30499.5  # L6: 30994 -> 30499.5
$ws.Cells.Item(6, 14).Value = -30725.5  # N6: -31220 -> -30725.5
$ws.Cells.Item(86, 8).Value = 1827.3334  # H86: 2493 -> 1827.3334
$ws.Cells.Item(86, 9).Value = 1827.3334  # I86: 2493 -> 1827.3334
$ws.Cells.Item(86, 11).Value = 1827.3334  # K86: 2493 -> 1827.3334
$ws.Cells.Item(86, 13).Value = -704.3334  # M86: -1370 -> -704.3334
$ws.Cells.Item(89, 8).Value = 1827.3334  # H89: 2493 -> 1827.3334
$ws.Cells.Item(89, 9).Value = 1827.3334  # I89: 2493 -> 1827.3334
$ws.Cells.Item(89, 11).Value = 9136.666999999999  # K89: 12465 -> 9136.666999999999
$ws.Cells.Item(89, 13).Value = -3520.666999999999  # M89: -6849 -> -3520.666999999999
$ws.Cells.Item(95, 8).Value = 20186.8  # H95: 21171 -> 20186.8
$ws.Cells.Item(95, 10).Value = 20186.8  # J95: 21171 -> 20186.8
$ws.Cells.Item(95, 12).Value = 20186.8  # L95: 21171 -> 20186.8
$ws.Cells.Item(95, 14).Value = -25678.8  # N95: -26663 -> -25678.8
$ws.Cells.Item(134, 8).Value = 2509.6428  # H134: 3327.111 -> 2509.6428
$ws.Cells.Item(134, 9).Value = 2433.4614  # I134: 3305.5 -> 2433.4614
$ws.Cells.Item(134, 11).Value = 7300.3842  # K134: 9916.5 -> 7300.3842
$ws.Cells.Item(134, 13).Value = -4765.3842  # M134: -7381.5 -> -4765.3842

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2158.0715  # H16: 2357.9285 -> 2158.0715
$ws.Cells.Item(16, 9).Value = 1704.125  # I16: 1919 -> 1704.125
$ws.Cells.Item(16, 10).Value = 2763.3333  # J16: 2796.8572 -> 2763.3333
$ws.Cells.Item(16, 11).Value = 1704.125  # K16: 1919 -> 1704.125
$ws.Cells.Item(16, 12).Value = 2763.3333  # L16: 2796.8572 -> 2763.3333
$ws.Cells.Item(16, 13).Value = -1417.125  # M16: -1632 -> -1417.125
$ws.Cells.Item(16, 14).Value = -3337.3333  # N16: -3370.8572 -> -3337.3333
$ws.Cells.Item(22, 8).Value = 32174.875  # H22: 36628.43 -> 32174.875
$ws.Cells.Item(22, 9).Value = 1180  # I22: 1225 -> 1180
$ws.Cells.Item(22, 11).Value = 1180  # K22: 1225 -> 1180
$ws.Cells.Item(22, 13).Value = -830  # M22: -875 -> -830
$ws.Cells.Item(31, 8).Value = 1724.8  # H31: 1863.625 -> 1724.8
$ws.Cells.Item(31, 9).Value = 1587  # I31: 1737.9231 -> 1587
$ws.Cells.Item(31, 10).Value = 2965  # J31: 2408.3333 -> 2965
$ws.Cells.Item(31, 11).Value = 1587  # K31: 1737.9231 -> 1587
$ws.Cells.Item(31, 12).Value = 2965  # L31: 2408.3333 -> 2965
$ws.Cells.Item(31, 13).Value = -1292  # M31: -1442.9231 -> -1292
$ws.Cells.Item(31, 14).Value = -3555  # N31: -2998.3333 -> -3555
$ws.Cells.Item(34, 8).Value = 1724.8  # H34: 1863.625 -> 1724.8
$ws.Cells.Item(34, 9).Value = 1587  # I34: 1737.9231 -> 1587
$ws.Cells.Item(34, 10).Value = 2965  # J34: 2408.3333 -> 2965
$ws.Cells.Item(34, 11).Value = 1587  # K34: 1737.9231 -> 1587
$ws.Cells.Item(34, 12).Value = 2965  # L34: 2408.3333 -> 2965
$ws.Cells.Item(34, 13).Value = -1385  # M34: -1535.9231 -> -1385
$ws.Cells.Item(34, 14).Value = -3369  # N34: -2812.3333 -> -3369
$ws.Cells.Item(58, 8).Value = 2873.1667  # H58: 3444 -> 2873.1667
$ws.Cells.Item(58, 9).Value = 2848  # I58: 3888.5 -> 2848
$ws.Cells.Item(58, 10).Value = 2999  # J58: 2999.5 -> 2999
$ws.Cells.Item(58, 11).Value = 2848  # K58: 3888.5 -> 2848
$ws.Cells.Item(58, 12).Value = 2999  # L58: 2999.5 -> 2999
$ws.Cells.Item(58, 13).Value = -2645  # M58: -3685.5 -> -2645
$ws.Cells.Item(58, 14).Value = -3405  # N58: -3405.5 -> -3405
$ws.Cells.Item(99, 8).Value = 3307.7778  # H99: 3754.7144 -> 3307.7778
$ws.Cells.Item(99, 9).Value = 1489.5  # I99: 1490 -> 1489.5
$ws.Cells.Item(99, 10).Value = 3827.2856  # J99: 4132.1665 -> 3827.2856
$ws.Cells.Item(99, 11).Value = 1489.5  # K99: 1490 -> 1489.5
$ws.Cells.Item(99, 12).Value = 3827.2856  # L99: 4132.1665 -> 3827.2856
$ws.Cells.Item(99, 13).Value = 8.5  # M99: 8 -> 8.5
$ws.Cells.Item(99, 14).Value = -6823.2856  # N99: -7128.1665 -> -6823.2856
$ws.Cells.Item(107, 8).Value = 858.6429000000001  # H107: 1008.7692 -> 858.6429000000001
$ws.Cells.Item(107, 9).Value = 703.7143  # I107: 750.5 -> 703.7143
$ws.Cells.Item(107, 10).Value = 1013.5714  # J107: 1230.1428 -> 1013.5714
$ws.Cells.Item(107, 11).Value = 703.7143  # K107: 750.5 -> 703.7143
$ws.Cells.Item(107, 12).Value = 1013.5714  # L107: 1230.1428 -> 1013.5714
$ws.Cells.Item(107, 13).Value = 1216.2857  # M107: 1169.5 -> 1216.2857
$ws.Cells.Item(107, 14).Value = -4853.5714  # N107: -5070.1428 -> -4853.5714
$ws.Cells.Item(113, 8).Value = 2158.0715  # H113: 2357.9285 -> 2158.0715
$ws.Cells.Item(113, 9).Value = 1704.125  # I113: 1919 -> 1704.125
$ws.Cells.Item(113, 10).Value = 2763.3333  # J113: 2796.8572 -> 2763.3333
$ws.Cells.Item(113, 11).Value = 1704.125  # K113: 1919 -> 1704.125
$ws.Cells.Item(113, 12).Value = 2763.3333  # L113: 2796.8572 -> 2763.3333
$ws.Cells.Item(113, 13).Value = 465.875  # M113: 251 -> 465.875
$ws.Cells.Item(113, 14).Value = -7103.3333  # N113: -7136.8572 -> -7103.3333
$ws.Cells.Item(122, 8).Value = 1513.8  # H122: 1757.6154 -> 1513.8
$ws.Cells.Item(122, 9).Value = 1120.875  # I122: 1285.1666 -> 1120.875
$ws.Cells.Item(122, 10).Value = 1962.8572  # J122: 2162.5715 -> 1962.8572
$ws.Cells.Item(122, 11).Value = 3362.625  # K122: 3855.4998 -> 3362.625
$ws.Cells.Item(122, 12).Value = 5888.571599999999  # L122: 6487.7145 -> 5888.571599999999
$ws.Cells.Item(122, 13).Value = -912.625  # M122: -1405.4998 -> -912.625
$ws.Cells.Item(122, 14).Value = -10788.5716  # N122: -11387.7145 -> -10788.5716
$ws.Cells.Item(126, 8).Value = 3307.7778  # H126: 3754.7144 -> 3307.7778
$ws.Cells.Item(126, 9).Value = 1489.5  # I126: 1490 -> 1489.5
$ws.Cells.Item(126, 10).Value = 3827.2856  # J126: 4132.1665 -> 3827.2856
$ws.Cells.Item(126, 11).Value = 4468.5  # K126: 4470 -> 4468.5
$ws.Cells.Item(126, 12).Value = 11481.8568  # L126: 12396.4995 -> 11481.8568
$ws.Cells.Item(126, 13).Value = -1998.5  # M126: -2000 -> -1998.5
$ws.Cells.Item(126, 14).Value = -16421.8568  # N126: -17336.4995 -> -16421.8568
$ws.Cells.Item(131, 8).Value = 69992  # H131: 99999 -> 69992
$ws.Cells.Item(131, 10).Value = 69992  # J131: 99999 -> 69992
$ws.Cells.Item(131, 12).Value = 69992  # L131: 99999 -> 69992
$ws.Cells.Item(131, 14).Value = -80072  # N131: -110079 -> -80072
$ws.Cells.Item(134, 8).Value = 4141.7144  # H134: 3955.5715 -> 4141.7144
$ws.Cells.Item(134, 9).Value = 3998.6667  # I134: 3955.5715 -> 3998.6667
$ws.Cells.Item(134, 10).Value = 5000  # J134: 0 -> 5000
$ws.Cells.Item(134, 11).Value = 11996.0001  # K134: 11866.7145 -> 11996.0001
$ws.Cells.Item(134, 12).Value = 15000  # L134: 0 -> 15000
$ws.Cells.Item(134, 13).Value = -9461.000100000001  # M134: -9331.7145 -> -9461.000100000001
$ws.Cells.Item(134, 14).Value = -20070  # N134: None -> -20070
$ws.Cells.Item(136, 8).Value = 2873.1667  # H136: 3444 -> 2873.1667
$ws.Cells.Item(136, 9).Value = 2848  # I136: 3888.5 -> 2848
$ws.Cells.Item(136, 10).Value = 2999  # J136: 2999.5 -> 2999
$ws.Cells.Item(136, 11).Value = 8544  # K136: 11665.5 -> 8544
$ws.Cells.Item(136, 12).Value = 8997  # L136: 8998.5 -> 8997
$ws.Cells.Item(136, 13).Value = -5994  # M136: -9115.5 -> -5994
$ws.Cells.Item(136, 14).Value = -14097  # N136: -14098.5 -> -14097

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 10235.277  # H121: 11473.5 -> 10235.277
$ws.Cells.Item(121, 9).Value = 15439.5  # I121: 20461.166 -> 15439.5
$ws.Cells.Item(121, 10).Value = 6071.9  # J121: 6080.9 -> 6071.9
$ws.Cells.Item(121, 11).Value = 46318.5  # K121: 61383.49800000001 -> 46318.5
$ws.Cells.Item(121, 12).Value = 18215.7  # L121: 18242.7 -> 18215.7
$ws.Cells.Item(121, 13).Value = -45008.5  # M121: -60073.49800000001 -> -45008.5
$ws.Cells.Item(121, 14).Value = -20835.7  # N121: -20862.7 -> -20835.7
$ws.Cells.Item(134, 8).Value = 14594.777  # H134: 12739.444 -> 14594.777
$ws.Cells.Item(134, 9).Value = 1525  # I134: 1466.3334 -> 1525
$ws.Cells.Item(134, 10).Value = 18329  # J134: 18376 -> 18329
$ws.Cells.Item(134, 11).Value = 4575  # K134: 4399.0002 -> 4575
$ws.Cells.Item(134, 12).Value = 54987  # L134: 55128 -> 54987
$ws.Cells.Item(134, 13).Value = 495  # M134: 670.9997999999996 -> 495
$ws.Cells.Item(134, 14).Value = -65127  # N134: -65268 -> -65127

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1299.5  # H80: 774.5 -> 1299.5
$ws.Cells.Item(80, 9).Value = 1299.5  # I80: 774.5 -> 1299.5
$ws.Cells.Item(80, 11).Value = 1299.5  # K80: 774.5 -> 1299.5
$ws.Cells.Item(80, 13).Value = -301.5  # M80: 223.5 -> -301.5
$ws.Cells.Item(83, 8).Value = 1299.5  # H83: 774.5 -> 1299.5
$ws.Cells.Item(83, 9).Value = 1299.5  # I83: 774.5 -> 1299.5
$ws.Cells.Item(83, 11).Value = 6497.5  # K83: 3872.5 -> 6497.5
$ws.Cells.Item(83, 13).Value = -1505.5  # M83: 1119.5 -> -1505.5
$ws.Cells.Item(122, 8).Value = 8251.25  # H122: 6762.2 -> 8251.25
$ws.Cells.Item(122, 9).Value = 6999  # I122: 4934.6665 -> 6999
$ws.Cells.Item(122, 11).Value = 20997  # K122: 14803.9995 -> 20997
$ws.Cells.Item(122, 13).Value = -18547  # M122: -12353.9995 -> -18547
$ws.Cells.Item(132, 8).Value = 1286.25  # H132: 1330 -> 1286.25
$ws.Cells.Item(132, 9).Value = 1270  # I132: 1318.3334 -> 1270
$ws.Cells.Item(132, 11).Value = 3810  # K132: 3955.0002 -> 3810
$ws.Cells.Item(132, 13).Value = -1280  # M132: -1425.0002 -> -1280

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6534.6665  # H7: 6816 -> 6534.6665
$ws.Cells.Item(7, 9).Value = 3288.7693  # I7: 3414.7273 -> 3288.7693
$ws.Cells.Item(7, 11).Value = 3288.7693  # K7: 3414.7273 -> 3288.7693
$ws.Cells.Item(7, 13).Value = -3176.7693  # M7: -3302.7273 -> -3176.7693
$ws.Cells.Item(16, 8).Value = 0  # H16: 5000 -> 0
$ws.Cells.Item(16, 10).Value = 0  # J16: 5000 -> 0
$ws.Cells.Item(16, 12).Value = 0  # L16: 5000 -> 0
$ws.Cells.Item(16, 14).ClearContents()  # N16 removed (was -5340)
$ws.Cells.Item(100, 8).Value = 4428.2856  # H100: 3994.7 -> 4428.2856
$ws.Cells.Item(100, 9).Value = 2499  # I100: 2512 -> 2499
$ws.Cells.Item(100, 10).Value = 5200  # J100: 4983.1665 -> 5200
$ws.Cells.Item(100, 11).Value = 2499  # K100: 2512 -> 2499
$ws.Cells.Item(100, 12).Value = 5200  # L100: 4983.1665 -> 5200
$ws.Cells.Item(100, 13).Value = -1958  # M100: -1971 -> -1958
$ws.Cells.Item(100, 14).Value = -6282  # N100: -6065.1665 -> -6282
$ws.Cells.Item(122, 8).Value = 6095.2085  # H122: 6103.5 -> 6095.2085
$ws.Cells.Item(122, 10).Value = 6579.3  # J122: 6599.2 -> 6579.3
$ws.Cells.Item(122, 12).Value = 19737.9  # L122: 19797.6 -> 19737.9
$ws.Cells.Item(122, 14).Value = -24637.9  # N122: -24697.6 -> -24637.9
$ws.Cells.Item(126, 8).Value = 6534.6665  # H126: 6816 -> 6534.6665
$ws.Cells.Item(126, 9).Value = 3288.7693  # I126: 3414.7273 -> 3288.7693
$ws.Cells.Item(126, 11).Value = 9866.3079  # K126: 10244.1819 -> 9866.3079
$ws.Cells.Item(126, 13).Value = -7396.3079  # M126: -7774.1819 -> -7396.3079
$ws.Cells.Item(132, 8).Value = 3634.4  # H132: 3906.125 -> 3634.4
$ws.Cells.Item(132, 9).Value = 2941.2856  # I132: 3098.8 -> 2941.2856
$ws.Cells.Item(132, 11).Value = 8823.856800000001  # K132: 9296.400000000001 -> 8823.856800000001
$ws.Cells.Item(132, 13).Value = -6293.856800000001  # M132: -6766.400000000001 -> -6293.856800000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 150231.67  # H119: 149999 -> 150231.67
$ws.Cells.Item(119, 10).Value = 150231.67  # J119: 149999 -> 150231.67
$ws.Cells.Item(119, 12).Value = 150231.67  # L119: 149999 -> 150231.67
$ws.Cells.Item(119, 14).Value = -159907.67  # N119: -159675 -> -159907.67
$ws.Cells.Item(126, 8).Value = 3991.1177  # H126: 4456.2666 -> 3991.1177
$ws.Cells.Item(126, 9).Value = 1759.3636  # I126: 2038.6666 -> 1759.3636
$ws.Cells.Item(126, 11).Value = 5278.0908  # K126: 6115.9998 -> 5278.0908
$ws.Cells.Item(126, 13).Value = -3278.0908  # M126: -3645.9998 -> -3278.0908

